$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 389
$ws.Range("I2").Value = 389
$ws.Range("K2").Value = 389
$ws.Range("M2").Value = -276
$ws.Range("H29").Value = 4245.923
$ws.Range("I29").Value = 991
$ws.Range("J29").Value = 6280.25
$ws.Range("K29").Value = 2973
$ws.Range("L29").Value = 18840.75
$ws.Range("M29").Value = -2692
$ws.Range("N29").Value = -19402.75
$ws.Range("H106").Value = 8770.5
$ws.Range("I106").Value = 3374.45
$ws.Range("J106").Value = 35750.75
$ws.Range("K106").Value = 3374.45
$ws.Range("L106").Value = 35750.75
$ws.Range("M106").Value = -2743.45
$ws.Range("N106").Value = -37012.75
$ws.Range("H137").Value = 3599.087
$ws.Range("I137").Value = 2962.75
$ws.Range("K137").Value = 8888.25
$ws.Range("M137").Value = -6338.25

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H43").Value = 30358
$ws.Range("J43").Value = 29404.75
$ws.Range("L43").Value = 29404.75
$ws.Range("N43").Value = -30030.75
$ws.Range("H61").Value = 5938.0967
$ws.Range("I61").Value = 5476.033
$ws.Range("K61").Value = 5476.033
$ws.Range("M61").Value = -5264.033
$ws.Range("H62").Value = 50000
$ws.Range("J62").Value = 50000
$ws.Range("L62").Value = 50000
$ws.Range("N62").Value = -51248
$ws.Range("H65").Value = 50000
$ws.Range("J65").Value = 50000
$ws.Range("L65").Value = 150000
$ws.Range("N65").Value = -156240
$ws.Range("H74").Value = 15874942
$ws.Range("I74").Value = 19609850
$ws.Range("K74").Value = 19609850
$ws.Range("M74").Value = -19608976
$ws.Range("H77").Value = 15874942
$ws.Range("I77").Value = 19609850
$ws.Range("K77").Value = 98049250
$ws.Range("M77").Value = -98044882
$ws.Range("H97").Value = 1008.2222
$ws.Range("I97").Value = 1376.909
$ws.Range("K97").Value = 1376.909
$ws.Range("M97").Value = -880.9090000000001
$ws.Range("H136").Value = 5938.0967
$ws.Range("I136").Value = 5476.033
$ws.Range("K136").Value = 16428.099
$ws.Range("M136").Value = -13878.099
$ws.Range("H138").Value = 80000
$ws.Range("J138").Value = 80000
$ws.Range("L138").Value = 80000
$ws.Range("N138").Value = -90280

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H10").Value = 2383.1667
$ws.Range("I10").Value = 2125
$ws.Range("J10").Value = 2899.5
$ws.Range("K10").Value = 2125
$ws.Range("L10").Value = 2899.5
$ws.Range("M10").Value = -1985
$ws.Range("N10").Value = -3179.5
$ws.Range("H134").Value = 2583.2666
$ws.Range("I134").Value = 1645.1666
$ws.Range("J134").Value = 6335.6665
$ws.Range("K134").Value = 4935.4998
$ws.Range("L134").Value = 19006.9995
$ws.Range("M134").Value = -2400.4998
$ws.Range("N134").Value = -24076.9995

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 3809.4
$ws.Range("I16").Value = 3365.6667
$ws.Range("J16").Value = 4475
$ws.Range("K16").Value = 3365.6667
$ws.Range("L16").Value = 4475
$ws.Range("M16").Value = -3078.6667
$ws.Range("N16").Value = -5049
$ws.Range("H58").Value = 3374.0625
$ws.Range("J58").Value = 8974.714
$ws.Range("L58").Value = 8974.714
$ws.Range("N58").Value = -9380.714
$ws.Range("H59").Value = 0
$ws.Range("J59").Value = 0
$ws.Range("L59").Value = 0
$ws.Range("N59").ClearContents()
$ws.Range("H113").Value = 3809.4
$ws.Range("I113").Value = 3365.6667
$ws.Range("J113").Value = 4475
$ws.Range("K113").Value = 3365.6667
$ws.Range("L113").Value = 4475
$ws.Range("M113").Value = -1195.6667
$ws.Range("N113").Value = -8815
$ws.Range("H136").Value = 3374.0625
$ws.Range("J136").Value = 8974.714
$ws.Range("L136").Value = 26924.142
$ws.Range("N136").Value = -32024.142

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H114").Value = 246.16667
$ws.Range("I114").Value = 245.8
$ws.Range("K114").Value = 737.4000000000001
$ws.Range("M114").Value = 2516.6
$ws.Range("H129").Value = 4632070.5
$ws.Range("I129").Value = 536.8182
$ws.Range("J129").Value = 11910195
$ws.Range("K129").Value = 1610.4546
$ws.Range("L129").Value = 35730585
$ws.Range("M129").Value = 3389.5454
$ws.Range("N129").Value = -35740585
$ws.Range("H133").Value = 8111.636
$ws.Range("I133").Value = 4006
$ws.Range("J133").Value = 11533
$ws.Range("K133").Value = 12018
$ws.Range("L133").Value = 34599
$ws.Range("M133").Value = -6958
$ws.Range("N133").Value = -44719
$ws.Range("H140").Value = 4022.2104
$ws.Range("I140").Value = 2745.2856
$ws.Range("J140").Value = 7597.6
$ws.Range("K140").Value = 8235.856800000001
$ws.Range("L140").Value = 22792.8
$ws.Range("M140").Value = -3055.856800000001
$ws.Range("N140").Value = -33152.8

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H52").Value = 0
$ws.Range("J52").Value = 0
$ws.Range("L52").Value = 0
$ws.Range("N52").ClearContents()
$ws.Range("H97").Value = 1303.3572
$ws.Range("I97").Value = 803.6667
$ws.Range("K97").Value = 803.6667
$ws.Range("M97").Value = -307.6667
$ws.Range("H102").Value = 3076.4285
$ws.Range("I102").Value = 2407.5
$ws.Range("J102").Value = 4748.75
$ws.Range("K102").Value = 2407.5
$ws.Range("L102").Value = 4748.75
$ws.Range("M102").Value = -785.5
$ws.Range("N102").Value = -7992.75
$ws.Range("H113").Value = 5042.8887
$ws.Range("I113").Value = 4641.3335
$ws.Range("J113").Value = 5846
$ws.Range("K113").Value = 4641.3335
$ws.Range("L113").Value = 5846
$ws.Range("M113").Value = -2471.3335
$ws.Range("N113").Value = -10186
$ws.Range("H122").Value = 15073.917
$ws.Range("I122").Value = 17175.445
$ws.Range("J122").Value = 8769.333000000001
$ws.Range("K122").Value = 51526.335
$ws.Range("L122").Value = 26307.999
$ws.Range("M122").Value = -49076.335
$ws.Range("N122").Value = -31207.999

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6626.3335
$ws.Range("I7").Value = 3405.1667
$ws.Range("K7").Value = 3405.1667
$ws.Range("M7").Value = -3293.1667
$ws.Range("H16").Value = 4662.8335
$ws.Range("I16").Value = 4662.8335
$ws.Range("K16").Value = 4662.8335
$ws.Range("M16").Value = -4492.8335
$ws.Range("H22").Value = 4233.091
$ws.Range("I22").Value = 1835.6
$ws.Range("K22").Value = 1835.6
$ws.Range("M22").Value = -1540.6
$ws.Range("H27").Value = 4233.091
$ws.Range("I27").Value = 1835.6
$ws.Range("K27").Value = 1835.6
$ws.Range("M27").Value = -1728.6
$ws.Range("H33").Value = 13756.125
$ws.Range("J33").Value = 13756.125
$ws.Range("L33").Value = 13756.125
$ws.Range("N33").Value = -14336.125
$ws.Range("H62").Value = 175166.67
$ws.Range("J62").Value = 175166.67
$ws.Range("L62").Value = 175166.67
$ws.Range("N62").Value = -176414.67
$ws.Range("H65").Value = 175166.67
$ws.Range("J65").Value = 175166.67
$ws.Range("L65").Value = 525500.01
$ws.Range("N65").Value = -531740.01
$ws.Range("H93").Value = 1482.1111
$ws.Range("I93").Value = 640
$ws.Range("J93").Value = 3166.3333
$ws.Range("K93").Value = 640
$ws.Range("L93").Value = 3166.3333
$ws.Range("M93").Value = 608
$ws.Range("N93").Value = -5662.3333
$ws.Range("H100").Value = 8849.933999999999
$ws.Range("I100").Value = 8993.777
$ws.Range("J100").Value = 8634.166999999999
$ws.Range("K100").Value = 8993.777
$ws.Range("L100").Value = 8634.166999999999
$ws.Range("M100").Value = -8452.777
$ws.Range("N100").Value = -9716.166999999999
$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("M122").ClearContents()
$ws.Range("H126").Value = 6626.3335
$ws.Range("I126").Value = 3405.1667
$ws.Range("K126").Value = 10215.5001
$ws.Range("M126").Value = -7745.500100000001
$ws.Range("H136").Value = 9450.0625
$ws.Range("I136").Value = 2854.7778
$ws.Range("J136").Value = 17929.715
$ws.Range("K136").Value = 8564.3334
$ws.Range("L136").Value = 53789.145
$ws.Range("M136").Value = -6014.3334
$ws.Range("N136").Value = -58889.145

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H5").Value = 16678666
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 16678666
$ws.Range("K5").Value = 0
$ws.Range("L5").Value = 16678666
$ws.Range("M5").ClearContents()
$ws.Range("N5").Value = -16678890
$ws.Range("H112").Value = 45672.668
$ws.Range("J112").Value = 45672.668
$ws.Range("L112").Value = 45672.668
$ws.Range("N112").Value = -48626.668
$ws.Range("H113").Value = 581.2
$ws.Range("J113").Value = 1200
$ws.Range("L113").Value = 3600
$ws.Range("N113").Value = -7940
$ws.Range("H115").Value = 50315.5
$ws.Range("J115").Value = 50315.5
$ws.Range("L115").Value = 50315.5
$ws.Range("N115").Value = -53449.5
$ws.Range("H116").Value = 220582.33
$ws.Range("J116").Value = 220582.33
$ws.Range("L116").Value = 220582.33
$ws.Range("N116").Value = -229760.33
$ws.Range("H132").Value = 7165.744
$ws.Range("I132").Value = 6352.925
$ws.Range("K132").Value = 19058.775
$ws.Range("M132").Value = -16528.775
$ws.Range("H136").Value = 9813.625
$ws.Range("I136").Value = 6400.8
$ws.Range("J136").Value = 15501.667
$ws.Range("K136").Value = 19202.4
$ws.Range("L136").Value = 46505.001
$ws.Range("M136").Value = -16652.4
$ws.Range("N136").Value = -51605.001
